$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above the current row 109 (the "note" row), pushing the
# note row down to row 110. Excel copies the formatting of the row above
# (row 108) onto the newly inserted row, which matches the target styles
# (s="3" for A, s="7" for B/C, s="8" for D/E).
$ws.Rows.Item(109).Insert()

# Populate the new row 109 with the day's data.
$ws.Range("A109").Value = 43964
$ws.Range("B109").Value = 297
$ws.Range("C109").Value = 36845
$ws.Range("D109").Value = 51
$ws.Range("E109").Value = 7488

# Update the print area defined name (localSheetId=0) so it covers the
# newly added row: $A$1:$E$110 -> $A$1:$E$111
$printAreaName = $null
foreach ($n in $wb.Names) {
    if ($n.Name -like "*Print_Area*") {
        $printAreaName = $n
    }
}
if ($printAreaName -eq $null) {
    $printAreaName = $wb.Names.Item(1)
}
$printAreaName.RefersTo = "=相談件数!`$A`$1:`$E`$111"

# Update the view selection on the (only) visible pane to A109, matching
# the post-edit sheet selection recorded in the workbook.
$ws.Activate()
$ws.Range("A109").Select()
